# GPLIM-3541: add Material Type as required header for Manifest uploads
#
# Adds a new "Material Type" column (G) to the manifest sheet: a styled
# header in G1 (bold white text on a black fill, centered - matching the
# look of the existing "SAMPLE_TYPE" header in F1) and "DNA:Genomic" filled
# down G2:G24 (centered, matching the style already used by column C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell G1: "Material Type" ------------------------------------
$ws.Range("G1").Value = "Material Type"
$ws.Range("G1").Font.Bold = $true
$ws.Range("G1").Font.Color = 16777215      # white
$ws.Range("G1").Interior.Color = 0         # black
$ws.Range("G1").HorizontalAlignment = -4108  # xlCenter

# --- Data cells G2:G24: "DNA:Genomic" ------------------------------------
$ws.Range("G2:G24").Value = "DNA:Genomic"
$ws.Range("G2:G24").HorizontalAlignment = -4108  # xlCenter

# --- Reflect the new column in the sheet's active selection -------------
$ws.Range("G1:G24").Select() | Out-Null
